$wb = $excel.ActiveWorkbook

# ALC row 6
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 1232.7273
$ws.Range("I6").Value = 330.9
$ws.Range("J6").Value = 10251
$ws.Range("K6").Value = 992.6999999999999
$ws.Range("L6").Value = 30753
$ws.Range("M6").Value = -880.6999999999999
$ws.Range("N6").Value = -30977

# ALC row 111
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H111").Value = 1668.0952
$ws.Range("I111").Value = 1733.1666
$ws.Range("J111").Value = 1581.3334
$ws.Range("K111").Value = 5199.4998
$ws.Range("L111").Value = 4744.0002
$ws.Range("M111").Value = -2132.4998
$ws.Range("N111").Value = -10878.0002

# ALC row 132
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 1202978.6
$ws.Range("I132").Value = 1302810.1
$ws.Range("K132").Value = 3908430.3
$ws.Range("M132").Value = -3905900.3

# ALC row 137
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 862.94116
$ws.Range("I137").Value = 769.3570999999999
$ws.Range("J137").Value = 1299.6666
$ws.Range("K137").Value = 2308.0713
$ws.Range("L137").Value = 3898.9998
$ws.Range("M137").Value = 241.9287000000004
$ws.Range("N137").Value = -8998.9998

# ARM row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6761.7344
$ws.Range("I32").Value = 2651.3076
$ws.Range("J32").Value = 25845.857
$ws.Range("K32").Value = 2651.3076
$ws.Range("L32").Value = 25845.857
$ws.Range("M32").Value = -2364.3076
$ws.Range("N32").Value = -26419.857

# ARM row 63
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 2224.5293
$ws.Range("J63").Value = 2427.1428
$ws.Range("L63").Value = 2427.1428
$ws.Range("N63").Value = -3799.1428

# ARM row 66
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H66").Value = 2224.5293
$ws.Range("J66").Value = 2427.1428
$ws.Range("L66").Value = 12135.714
$ws.Range("N66").Value = -18999.714

# ARM row 97
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 852
$ws.Range("I97").Value = 788.0625
$ws.Range("J97").Value = 1193
$ws.Range("K97").Value = 788.0625
$ws.Range("L97").Value = 1193
$ws.Range("M97").Value = -292.0625
$ws.Range("N97").Value = -2185

# CRP row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1576.909
$ws.Range("I31").Value = 1066.6333
$ws.Range("K31").Value = 1066.6333
$ws.Range("M31").Value = -771.6333

# CRP row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 1576.909
$ws.Range("I34").Value = 1066.6333
$ws.Range("K34").Value = 1066.6333
$ws.Range("M34").Value = -864.6333

# CRP row 74
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H74").Value = 76934970
$ws.Range("J74").Value = 83346024
$ws.Range("L74").Value = 83346024
$ws.Range("N74").Value = -83347772

# CRP row 77
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H77").Value = 76934970
$ws.Range("J77").Value = 83346024
$ws.Range("L77").Value = 250038072
$ws.Range("N77").Value = -250046808

# CRP row 88
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H88").Value = 35171.5
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 35171.5
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 35171.5
$ws.Range("N88").Value = -35983.5
$ws.Range("M88").ClearContents()

# CRP row 91
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H91").Value = 35171.5
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 35171.5
$ws.Range("K91").Value = 0
$ws.Range("L91").Value = 35171.5
$ws.Range("N91").Value = -37979.5
$ws.Range("M91").ClearContents()

# CRP row 132
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 4891.129
$ws.Range("I132").Value = 5120.269
$ws.Range("J132").Value = 3699.6
$ws.Range("K132").Value = 15360.807
$ws.Range("L132").Value = 11098.8
$ws.Range("M132").Value = -12830.807
$ws.Range("N132").Value = -16158.8

# CRP row 134
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 1702.9667
$ws.Range("I134").Value = 1699.6207
$ws.Range("J134").Value = 1800
$ws.Range("K134").Value = 5098.8621
$ws.Range("L134").Value = 5400
$ws.Range("M134").Value = -2563.8621
$ws.Range("N134").Value = -10470

# CUL row 120
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H120").Value = 15124.875
$ws.Range("I120").Value = 9999.5
$ws.Range("K120").Value = 29998.5
$ws.Range("M120").Value = -25160.5

# CUL row 125
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H125").Value = 3922
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 3922
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = 11766
$ws.Range("N125").Value = -21606
$ws.Range("M125").ClearContents()

# GSM row 70
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6585.933
$ws.Range("I70").Value = 5872.5
$ws.Range("K70").Value = 5872.5
$ws.Range("M70").Value = -5602.5

# GSM row 73
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 6585.933
$ws.Range("I73").Value = 5872.5
$ws.Range("K73").Value = 5872.5
$ws.Range("M73").Value = -4936.5

# GSM row 80
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3021
$ws.Range("I80").Value = 3568.3333
$ws.Range("J80").Value = 2200
$ws.Range("K80").Value = 3568.3333
$ws.Range("L80").Value = 2200
$ws.Range("M80").Value = -2570.3333
$ws.Range("N80").Value = -4196

# GSM row 83
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 3021
$ws.Range("I83").Value = 3568.3333
$ws.Range("J83").Value = 2200
$ws.Range("K83").Value = 17841.6665
$ws.Range("L83").Value = 11000
$ws.Range("M83").Value = -12849.6665
$ws.Range("N83").Value = -20984

# GSM row 132
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 4440.45
$ws.Range("I132").Value = 5184
$ws.Range("J132").Value = 2705.5
$ws.Range("K132").Value = 15552
$ws.Range("L132").Value = 8116.5
$ws.Range("M132").Value = -13022
$ws.Range("N132").Value = -13176.5

# LTW row 22
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 690
$ws.Range("I22").Value = 650
$ws.Range("J22").Value = 750
$ws.Range("K22").Value = 650
$ws.Range("L22").Value = 750
$ws.Range("M22").Value = -355
$ws.Range("N22").Value = -1340

# LTW row 27
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 690
$ws.Range("I27").Value = 650
$ws.Range("J27").Value = 750
$ws.Range("K27").Value = 650
$ws.Range("L27").Value = 750
$ws.Range("M27").Value = -543
$ws.Range("N27").Value = -964

# LTW row 61
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1619.0476
$ws.Range("I61").Value = 1309.3334
$ws.Range("J61").Value = 2032
$ws.Range("K61").Value = 1309.3334
$ws.Range("L61").Value = 2032
$ws.Range("M61").Value = -1107.3334
$ws.Range("N61").Value = -2436

# LTW row 113
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 1619.0476
$ws.Range("I113").Value = 1309.3334
$ws.Range("J113").Value = 2032
$ws.Range("K113").Value = 1309.3334
$ws.Range("L113").Value = 2032
$ws.Range("M113").Value = 860.6666
$ws.Range("N113").Value = -6372

# LTW row 120
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H120").Value = 0
$ws.Range("J120").Value = 0
$ws.Range("L120").Value = 0
$ws.Range("N120").ClearContents()
